$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Add a new "description" column to the problem table (Table4)
$lo = $ws3.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

# Set header + new column width
$ws3.Range("D1").Value = "description"
$ws3.Columns.Item(4).ColumnWidth = 30

# Fill in description values for each row
$ws3.Range("D2").Value = "maximization of total profit"
$ws3.Range("D3").Value = "energy use less than endowment"
$ws3.Range("D4").Value = "positive products supply"
$ws3.Range("D5").Value = "energy use per unit of product"

# Update the selected cell on the problem sheet
$ws3.Activate()
$ws3.Range("C8").Select()
